$p = $ppt.ActivePresentation

# --- Locate the shape containing the GitHub repository URL (slide 2) ---
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$oldUrl = "https://github.com/sraorao/MSD_R_course_TT2022"
$newUrl = "https://github.com/sraorao/MSD_R_course_June2022"

# --- Step 1: update the URL text in place (keeps run formatting + hyperlink) ---
$fullText = $tr.Text
$zeroBasedIdx = $fullText.IndexOf($oldUrl)
if ($zeroBasedIdx -lt 0) {
    throw "Could not find the expected URL text in the shape"
}
$urlStart = $zeroBasedIdx + 1
$urlChars = $tr.Characters($urlStart, $oldUrl.Length)
$urlChars.Text = $newUrl

# --- Step 2: append a trailing space right after the (now renamed) URL run ---
# Re-resolve the paragraph that holds the URL run, then InsertAfter appends using
# that paragraph's last run formatting (size/colour/font), which matches the URL run.
$fullText2 = $tr.Text
$zeroBasedIdx2 = $fullText2.IndexOf($newUrl)
$paraCount = 0
$targetPara = $null
for ($i = 1; $i -le 50; $i++) {
    $thisPara = $tr.Paragraphs($i, 1)
    if ($thisPara -eq $null) { break }
    if ($thisPara.Text.IndexOf($newUrl) -ge 0) {
        $targetPara = $thisPara
        $paraCount = $i
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find the paragraph holding the renamed URL"
}

$spaceRange = $targetPara.InsertAfter(" ")

# --- Step 3: strip the hyperlink from just the newly appended trailing space ---
$refreshedPara = $tr.Paragraphs($paraCount, 1)
$lastChar = $refreshedPara.Characters($refreshedPara.Length - 1, 1)
$lastChar.ActionSettings.Item(1).Hyperlink.Address = ""
